$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 11

# Copy formatting from the row above (A10) so the new year label picks up
# the same bold/bordered/centered style used for every other year cell.
$ws.Range("A10").Copy()
$ws.Range("A$row").PasteSpecial(-4122)
$ws.Range("A$row").Value = "2021年"
$ws.Range("B$row").Value = 13608.68
$ws.Range("C$row").Value = 3547.44
$ws.Range("D$row").Value = 535.33

# This industry has no reported figure for 2021 — write it as an explicit
# empty text value (rather than leaving the cell untouched) and then reset
# the style to Normal so it doesn't pick up a stray quote-prefix format.
$ws.Range("E$row").Value = "'"
$ws.Range("E$row").Style = "Normal"

$ws.Range("F$row").Value = 10571.53
$ws.Range("G$row").Value = 13679.97
$ws.Range("H$row").Value = 2029.01
$ws.Range("I$row").Value = 7061.48
$ws.Range("J$row").Value = 1993.75
$ws.Range("K$row").Value = 2473.69
$ws.Range("L$row").Value = 1608.89
$ws.Range("M$row").Value = 179.2
$ws.Range("N$row").Value = 3353.88
$ws.Range("O$row").Value = 8227.690000000001
$ws.Range("P$row").Value = 616.54
$ws.Range("Q$row").Value = 2541.92
$ws.Range("R$row").Value = 8349.32
$ws.Range("S$row").Value = 353.59
$ws.Range("T$row").Value = 10264.69
$ws.Range("U$row").Value = 33.18
$ws.Range("V$row").Value = 4495.02
$ws.Range("W$row").Value = 778.0700000000001
$ws.Range("X$row").Value = 4616.05
$ws.Range("Y$row").Value = 25963.55
$ws.Range("Z$row").Value = 2167.86
$ws.Range("AA$row").Value = 6189.85
$ws.Range("AB$row").Value = 19.58
$ws.Range("AC$row").Value = 249374.7
$ws.Range("AD$row").Value = 7562.59
$ws.Range("AE$row").Value = 3541.1
$ws.Range("AF$row").Value = 21816.21
$ws.Range("AG$row").Value = 14787.72
$ws.Range("AH$row").Value = 2870.82
$ws.Range("AI$row").Value = 2498.37
$ws.Range("AJ$row").Value = 217.58
$ws.Range("AK$row").Value = 15201.87
$ws.Range("AL$row").Value = 3216.34
$ws.Range("AM$row").Value = 22940.79
$ws.Range("AN$row").Value = 1080.72
$ws.Range("AO$row").Value = 3610.64
$ws.Range("AP$row").Value = 13341.84
$ws.Range("AQ$row").Value = 1427.96
